# Regenerate the "K" column (strikeouts, column G) values for the save_data
# sheet. This mirrors a re-run of the upstream data pipeline that now
# computes K (strikeouts) straight from the box score instead of the old
# "Strike#" derived figure, so the raw integers below are the freshly
# calculated values to be written back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for G2:G59, in row order.
$newK = @(
    2,1,1,2,0,2,0,1,0,3,
    2,0,2,0,2,1,1,1,1,0,
    3,1,1,1,3,2,1,0,0,2,
    3,3,1,0,1,1,1,0,1,1,
    0,0,1,0,1,1,2,2,1,2,
    1,1,1,1,2,0,1,1
)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
